$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.308.58'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.887.83'
$ws.Range("E3").Value = '  -1.28%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.72'
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4682'
$ws.Range("E7").Value = '  -2.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06594'
$ws.Range("E9").Value = '  -1.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.78'
$ws.Range("E10").Value = '  +11.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07777'
$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("E12").Value = '  -3.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.883.57'
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.099'
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6760'
$ws.Range("E15").Value = '  +1.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.35'
$ws.Range("E16").Value = '  +11.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.318.03'
$ws.Range("E17").Value = '  -0.77%  '

$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.125.40'
$ws.Range("E20").Value = '  -1.54%  '

$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007292'
$ws.Range("E22").Value = '  -2.32%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.185'
$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.415'
$ws.Range("E25").Value = '  +1.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.75'
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.26'
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.986'
$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09759'
$ws.Range("E30").Value = '  -3.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.372'
$ws.Range("E31").Value = '  -7.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.484'
$ws.Range("E32").Value = '  -1.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.139'
$ws.Range("E33").Value = '  -2.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04672'
$ws.Range("E34").Value = '  -0.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7074'
$ws.Range("E35").Value = '  -2.54%  '

$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("E38").Value = '  -2.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.716'
$ws.Range("E39").Value = '  +7.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.532'
$ws.Range("E40").Value = '  -2.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.21'
$ws.Range("E41").Value = '  -3.16%  '

$ws.Range("E42").Value = '  +1.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.965'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.72'
$ws.Range("E45").Value = '  -1.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4184'
$ws.Range("E46").Value = '  -1.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '992.13'
$ws.Range("E47").Value = '  +8.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.284'
$ws.Range("E48").Value = '  -0.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.252'
$ws.Range("E49").Value = '  +6.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.96'
$ws.Range("E50").Value = '  -2.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1150'
$ws.Range("E51").Value = '  -3.87%  '
